# Applies the "output generated" refresh: a handful of "想去人数" (interest
# count) cells in column F tick up by 1 (or a few points) on three of the
# four sheets, matching the scraped gh-pages rebuild diff.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 ("Exhibitions")
$wsExpo = $wb.Worksheets.Item(1)
$wsExpo.Range("F10").Value = 1662
$wsExpo.Range("F11").Value = 1662
$wsExpo.Range("F18").Value = 608
$wsExpo.Range("F21").Value = 7347
$wsExpo.Range("F22").Value = 8178
$wsExpo.Range("F41").Value = 753
$wsExpo.Range("F44").Value = 354

# Sheet 2: 演出 ("Performances")
$wsShow = $wb.Worksheets.Item(2)
$wsShow.Range("F3").Value = 34

# Sheet 4: 全部类型 ("All types") -- union of the other sheets, updated in
# lockstep with the same underlying events.
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F8").Value = 34
$wsAll.Range("F12").Value = 1662
$wsAll.Range("F13").Value = 1662
$wsAll.Range("F19").Value = 608
$wsAll.Range("F24").Value = 7347
$wsAll.Range("F25").Value = 8178
$wsAll.Range("F39").Value = 753
$wsAll.Range("F44").Value = 354
